$wb = $excel.ActiveWorkbook

# --- Update the "Source" sheet CSV path references ---
$ws2 = $wb.Worksheets.Item("Source")
$ws2.Range("C2").Value2 = "examples/tutorial/operator-lines.csv"
$ws2.Range("C4").Value2 = "examples/tutorial/operator-lines.csv"
$ws2.Range("C6").Value2 = "examples/tutorial/operator-lines.csv"
$ws2.Range("C8").Value2 = "examples/tutorial/routes.csv"

# --- Change the active sheet / selection so that "Source" becomes the
#     active tab (previously "Function" was active) ---
$ws2.Activate()
$ws2.Range("B11").Select()
